$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a data-entry error: Greg's Pick for the Colts @ Titans game (row 5)
# was mistakenly carried over as the previous row's team; it should be the
# (correctly predicted) away team, the Indianapolis Colts.
$ws.Range("G5").Value = "Indianapolis Colts"

# --- Highlight (yellow fill) the cells that correctly predicted the away
# team for each finished game, across the model columns (C:G). These cells
# keep their existing thin right border (row 2-16 style), so re-applying the
# interior color is all that's needed - Excel will mint a new cell style that
# reuses that border with the new fill.
$highlightNormal = @("C3:G3", "C4:G4", "C6:G6", "D9:F9", "C10", "C11", "E11:F11", "E14:G14", "C15:F15", "D16")
$rngNormal = $ws.Range($highlightNormal[0])
for ($i = 1; $i -lt $highlightNormal.Length; $i++) {
    $rngNormal = $excel.Union($rngNormal, $ws.Range($highlightNormal[$i]))
}
$rngNormal.Interior.Color = 65535

# Same highlight treatment for the last data row (row 17), which carries a
# different (bottom+right) border style.
$highlightLast = @("C17:E17", "G17")
$rngLast = $ws.Range($highlightLast[0])
for ($i = 1; $i -lt $highlightLast.Length; $i++) {
    $rngLast = $excel.Union($rngLast, $ws.Range($highlightLast[$i]))
}
$rngLast.Interior.Color = 65535

# --- Add a new summary row (18) with each model's pick accuracy so far,
# formatted as a percentage.
$ws.Range("C18:G18").NumberFormat = "0.00%"
$ws.Range("C18").Value = 0.5625
$ws.Range("D18").Value = 0.5625
$ws.Range("E18").Value = 0.5
$ws.Range("F18").Value = 0.5625
$ws.Range("G18").Value = 0.6875

# --- Restore the selection to where the author left off.
$ws.Range("I20").Select()
